$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect the new "through" date
$ws.Name = "Through 2022-10-19"

# Update header label in I1 to match new date
$ws.Range("I1").Value = "2022 (through 10-19)"

# Update November (row 11) and Total (row 14) values for column I
$ws.Range("I11").Value = 64
$ws.Range("I14").Value = 1341
